# Refresh the Marketo usage data pulled into Sheet1 (log pushed 23/9/22).
#
# - The "Campaign Data" block (rows 8-14), which was still full of empty
#   placeholder rows, now carries the real campaign counts/labels.
# - Row 15 loses its stray "Total" label.
# - Rows 16-17 pick up "Segmentations"/"Leads" labels (Database Data block).
# - The old "Program Data" block (rows 18-20) collapses down to just a
#   couple of residual cells as the source report layout shifted up.
# - A few totals further down the sheet were refreshed, and the "Library"
#   usage row became "Programe Library" reporting text "False".
# - Selection moves to A29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $text) {
    # Leading apostrophe forces a literal text entry (stops Excel from
    # re-interpreting numeric-looking strings like "8" or "False" as a
    # number/boolean).
    $ws.Range($addr).Value = "'" + $text
}

function Set-EmptyText($addr) {
    $ws.Range($addr).Value = "'"
}

# --- Campaign Data block (rows 8-14): fill in the real figures ---
Set-Text "A8" "Campaign Data"
Set-Text "C8" "Default"
Set-Text "D8" "Automation"
Set-Text "E8" "WP_DEMO"

Set-Text "A9" "All Triggered Campaigns"
$ws.Range("B9").Value = 8
Set-Text "C9" "8"
Set-Text "D9" "0"
Set-Text "E9" "0"

Set-Text "A10" "Active Triggered Campaigns"
$ws.Range("B10").Value = 1
Set-Text "C10" "1"
Set-Text "D10" "0"
Set-Text "E10" "0"

Set-Text "A11" "Batch Campaigns - Repeating Schedule"
Set-Text "C11" "0"
Set-Text "D11" "0"
Set-Text "E11" "0"

Set-Text "A12" "All Batch Campaigns"
$ws.Range("B12").Value = 2
Set-Text "C12" "1"
Set-Text "D12" "1"
Set-Text "E12" "0"

Set-Text "A13" "All Campaigns"
$ws.Range("B13").Value = 10
Set-Text "C13" "9"
Set-Text "D13" "1"
Set-Text "E13" "0"

Set-Text "A14" "Active Campaigns"
$ws.Range("B14").Value = 1
Set-Text "C14" "1"
Set-Text "D14" "0"
Set-Text "E14" "0"

# --- Row 15: the "Total" label in B15 no longer applies here ---
$ws.Range("B15").ClearContents()

# --- Database Data block (rows 16-17): Segmentations / Leads ---
Set-Text "A16" "Segmentations"
$ws.Range("B16").ClearContents()
$ws.Range("C16").Value = 2

Set-Text "A17" "Leads"
$ws.Range("B17").ClearContents()

# --- Old "Program Data" block (rows 18-20) collapses down ---
$ws.Range("A18").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("D18").ClearContents()
Set-EmptyText "E18"

$ws.Range("A19").ClearContents()
$ws.Range("B19").Value = 0
$ws.Range("C19").ClearContents()
$ws.Range("D19").ClearContents()
Set-EmptyText "E19"

$ws.Range("A20").ClearContents()
Set-EmptyText "B20"

# --- Refreshed totals further down the sheet ---
$ws.Range("B27").Value = 3
$ws.Range("B28").Value = 2
$ws.Range("B29").Value = 3

# Row 30 used to be the "Library" usage count; it is now "Programe Library"
# reported as text "False" rather than a number.
Set-Text "A30" "Programe Library"
Set-Text "B30" "False"

# --- Move the saved selection to A29, matching where the author left off ---
$ws.Range("A29").Select() | Out-Null
